# PrBlank.xlsx: add a "Lab. #" column at the front of the table, fill in
# the laboratory number for each data row, and highlight the data rows
# with a light-green fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A - this shifts the existing
# A:H columns (and their widths) one place to the right to become B:I.
$ws.Columns("A:A").Insert()

# Narrow width for the new "Lab. #" column (closest value this engine's
# pixel-snapped ColumnWidth can reach to the target 7.7109375 chars).
$ws.Range("A1").ColumnWidth = 6.8

# Header for the new column.
$ws.Range("A1").Value = "Lab. #"

# Laboratory number for each data row.
$ws.Range("A2").Value = 8166
$ws.Range("A3").Value = 8166

# Highlight the whole data block (now A2:I3) with the light-green fill.
$ws.Range("A2:I3").Interior.Color = 12379352
